$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Formula = "65.461.99"
$ws.Range("E2").Formula = "  +2.03%  "

# Row 3
$ws.Range("D3").Formula = "3.172.46"
$ws.Range("E3").Formula = "  +3.63%  "

# Row 4
$ws.Range("E4").Formula = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "  +1.14%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "  +7.34%  "

# Row 7
$ws.Range("E7").Formula = "  -0.12%  "

# Row 8
$ws.Range("D8").Formula = "3.164.13"
$ws.Range("E8").Formula = "  +3.76%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.504"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Formula = "  +3.08%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Formula = "  +15.64%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.162"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Formula = "  +3.03%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.473"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "  +2.54%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Formula = "  +3.99%  "

# Row 14
$ws.Range("E14").Formula = "  +2.11%  "

# Row 15
$ws.Range("D15").Formula = "3.674.91"
$ws.Range("E15").Formula = "  +3.59%  "

# Row 16
$ws.Range("D16").Formula = "65.475.36"
$ws.Range("E16").Formula = "  +1.97%  "

# Row 17
$ws.Range("E17").Formula = "  +10.17%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.114"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Formula = "  +2.40%  "

# Row 19
$ws.Range("D19").Formula = "3.166.44"
$ws.Range("E19").Formula = "  +3.26%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "  +4.45%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "  +3.58%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.720"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Formula = "  +5.68%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "  +5.70%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "  +5.03%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Formula = "  +2.65%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Formula = "  +0.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Formula = "  +18.65%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Formula = "  +3.08%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Formula = "  +6.61%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Formula = "  -0.22%  "

# Row 31
$ws.Range("B31").Formula = "EthereumClassic"
$ws.Range("C31").Formula = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Formula = "  +1.92%  "

# Row 32
$ws.Range("B32").Formula = "Stacks"
$ws.Range("C32").Formula = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Formula = "  +3.60%  "

# Row 33
$ws.Range("E33").Formula = "  +4.63%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "555.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Formula = "  +9.94%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Formula = "  +2.50%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Formula = "  +6.43%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0456"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Formula = "  +12.97%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "53.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Formula = "  +1.04%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0839"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Formula = "  +6.35%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "  +13.15%  "

# Row 41
$ws.Range("E41").Formula = "  +3.77%  "

# Row 42
$ws.Range("D42").Formula = "3.076.09"
$ws.Range("E42").Formula = "  +5.84%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "  +1.94%  "

# Row 44
$ws.Range("B44").Formula = "Fetch.AI"
$ws.Range("C44").Formula = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Formula = "  +12.49%  "

# Row 45
$ws.Range("B45").Formula = "TheGraph"
$ws.Range("C45").Formula = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.267"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "  +10.13%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "  +4.58%  "

# Row 48
$ws.Range("D48").Formula = "0.0₃0534"
$ws.Range("E48").Formula = "  -0.32%  "

# Row 49
$ws.Range("E49").Formula = "  +3.53%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "120.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Formula = "  -0.61%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Formula = "  +5.16%  "
